$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OCT-2020")

# Order of writes controls the new shared-string index allocation, so it is
# chosen to match the target: "Hayaai " (48), "To create new table..." (49),
# "To add datas..." (50), "Analyse the Token..." (51), "Analyse and Create..." (52)
$ws.Range("C9").Value = "Hayaai "

$ws.Range("D10").Value = "To create new table and insert datas as like test db in new instance db."
$ws.Range("D9").Value = "To add datas in the testdb and new instance db and support Haayai app for login screen development"

$ws.Range("C10").Value = "Hayaai "
$ws.Range("C11").Value = "Hayaai "
$ws.Range("D11").Value = "Analyse the Token system requirements and checked in GSS"

$ws.Range("C12").Value = "Hayaai "
$ws.Range("D12").Value = "Analyse and Create a table for Token system for Hayaai app."

$ws.Range("D12").Select()
